$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ([string]::IsNullOrEmpty($val)) { continue }

    $parts = $val -split ','
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $sysParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq 'system') {
            $sysParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($sysParts.Count -gt 0 -and $otherParts.Count -gt 0) {
        $newParts = $sysParts + $otherParts
        $newVal = [string]::Join(', ', $newParts)
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
